$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add()

$ws.Range("B5").Value = "~ for each contact generate a filled copy of the envelop for each address set`n~ for each contact generate carbon copies of the letter if the contact has more than one address set"
$ws.Range("A5").Value = "[1.6]"
$ws.Range("C5").Value = 43248

$ws.Range("B5").WrapText = $true
$ws.Range("C5").NumberFormat = $ws.Range("C4").NumberFormat

$ws.Rows.Item(5).RowHeight = 45

$ws.Range("A5:C5").Select()
